# Weekly update: insert two new price rows for Femacal de La Calera - Alcachofa
# (Argentina(o)/Primera and Española/Primera) above the previous first data
# block, shifting all subsequent rows down by two and extending the used
# range from A1:R564 to A1:R566.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 491; everything below (old rows 491-564)
# moves down to 493-566.
$ws.Rows("491:492").Insert()

# New row 491: Argentina(o) / Primera
$ws.Range("A491").Value = 3
$ws.Range("B491").Value = "Femacal de La Calera"
$ws.Range("C491").Value = "Coquimbo"
$ws.Range("D491").Value = 45127
$ws.Range("E491").Value = 5
$ws.Range("F491").Value = 100112013
$ws.Range("G491").Value = "Alcachofa"
$ws.Range("H491").Value = "Argentina(o)"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 115
$ws.Range("K491").Value = 14000
$ws.Range("L491").Value = 14500
$ws.Range("M491").Value = 14217
$ws.Range("N491").Value = "$/caja 50 unidades"
$ws.Range("O491").Value = "Provincia de Limarí"
$ws.Range("P491").Value = 284
$ws.Range("Q491").Value = 50
$ws.Range("R491").Value = "Hortaliza"

# New row 492: Española / Primera
$ws.Range("A492").Value = 3
$ws.Range("B492").Value = "Femacal de La Calera"
$ws.Range("C492").Value = "Coquimbo"
$ws.Range("D492").Value = 45127
$ws.Range("E492").Value = 5
$ws.Range("F492").Value = 100112013
$ws.Range("G492").Value = "Alcachofa"
$ws.Range("H492").Value = "Española"
$ws.Range("I492").Value = "Primera"
$ws.Range("J492").Value = 115
$ws.Range("K492").Value = 15000
$ws.Range("L492").Value = 15500
$ws.Range("M492").Value = 15217
$ws.Range("N492").Value = "$/caja 30 unidades"
$ws.Range("O492").Value = "Provincia de Limarí"
$ws.Range("P492").Value = 507
$ws.Range("Q492").Value = 30
$ws.Range("R492").Value = "Hortaliza"
